$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-18 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-19 Monday", 2) | Out-Null
$d.Content.Find.Execute("59÷9=6, 5", $true, $false, $false, $false, $false, $true, 1, $false, "76÷5=15, 1", 2) | Out-Null
$d.Content.Find.Execute("39÷4=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "26÷7=3, 5", 2) | Out-Null
$d.Content.Find.Execute("89÷7=12, 5", $true, $false, $false, $false, $false, $true, 1, $false, "26÷5=5, 1", 2) | Out-Null
$d.Content.Find.Execute("19÷5=3, 4", $true, $false, $false, $false, $false, $true, 1, $false, "65÷5=13, 0", 2) | Out-Null
$d.Content.Find.Execute("21÷6=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "20÷7=2, 6", 2) | Out-Null
$d.Content.Find.Execute("63÷4=15, 3", $true, $false, $false, $false, $false, $true, 1, $false, "12÷4=3, 0", 2) | Out-Null
$d.Content.Find.Execute("78÷5=15, 3", $true, $false, $false, $false, $false, $true, 1, $false, "79÷4=19, 3", 2) | Out-Null
$d.Content.Find.Execute("61÷9=6, 7", $true, $false, $false, $false, $false, $true, 1, $false, "97÷8=12, 1", 2) | Out-Null
$d.Content.Find.Execute("57÷5=11, 2", $true, $false, $false, $false, $false, $true, 1, $false, "21÷8=2, 5", 2) | Out-Null
$d.Content.Find.Execute("16÷2=8, 0", $true, $false, $false, $false, $false, $true, 1, $false, "90÷3=30, 0", 2) | Out-Null
$d.Content.Find.Execute("76÷2=38, 0", $true, $false, $false, $false, $false, $true, 1, $false, "74÷8=9, 2", 2) | Out-Null
$d.Content.Find.Execute("34÷2=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "72÷2=36, 0", 2) | Out-Null
$d.Content.Find.Execute("33÷2=16, 1", $true, $false, $false, $false, $false, $true, 1, $false, "46÷3=15, 1", 2) | Out-Null
$d.Content.Find.Execute("93÷6=15, 3", $true, $false, $false, $false, $false, $true, 1, $false, "63÷3=21, 0", 2) | Out-Null
$d.Content.Find.Execute("83÷8=10, 3", $true, $false, $false, $false, $false, $true, 1, $false, "33÷6=5, 3", 2) | Out-Null
$d.Content.Find.Execute("28÷5=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "69÷4=17, 1", 2) | Out-Null
$d.Content.Find.Execute("96÷8=12, 0", $true, $false, $false, $false, $false, $true, 1, $false, "75÷3=25, 0", 2) | Out-Null
$d.Content.Find.Execute("75÷4=18, 3", $true, $false, $false, $false, $false, $true, 1, $false, "98÷7=14, 0", 2) | Out-Null
$d.Content.Find.Execute("34÷7=4, 6", $true, $false, $false, $false, $false, $true, 1, $false, "20÷3=6, 2", 2) | Out-Null
$d.Content.Find.Execute("90÷2=45, 0", $true, $false, $false, $false, $false, $true, 1, $false, "24÷7=3, 3", 2) | Out-Null
$d.Content.Find.Execute("83÷3=27, 2", $true, $false, $false, $false, $false, $true, 1, $false, "32÷8=4, 0", 2) | Out-Null
$d.Content.Find.Execute("87÷9=9, 6", $true, $false, $false, $false, $false, $true, 1, $false, "41÷6=6, 5", 2) | Out-Null
$d.Content.Find.Execute("22÷2=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "95÷7=13, 4", 2) | Out-Null
$d.Content.Find.Execute("35÷5=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "96÷7=13, 5", 2) | Out-Null
$d.Content.Find.Execute("14÷3=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "82÷3=27, 1", 2) | Out-Null
